$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -0.000003584038161252856
$ws.Range("D2").Value = -0.0007063712514252529
$ws.Range("F2").Value = -0.000003584038161252856
$ws.Range("G2").Value = -0.0000179201908063753
$ws.Range("C3").Value = -0.00005427818574088583
$ws.Range("D3").Value = -0.0002171127429635433
$ws.Range("E3").Value = -0.00006636440568952295
$ws.Range("F3").Value = -0.0001006096718163541
$ws.Range("G3").Value = -0.0002713909287042071
$ws.Range("C4").Value = -0.00000345765962791577
$ws.Range("D4").Value = -0.00001383063851166308
$ws.Range("E4").Value = -0.01581441097099656
$ws.Range("F4").Value = 0.0158762674836419
$ws.Range("G4").Value = -0.00001728829813957972
$ws.Range("C5").Value = -0.0007098437176864536
$ws.Range("D5").Value = -0.000001192260392191002
$ws.Range("F5").Value = -0.0000002980650980477506
$ws.Range("G5").Value = -0.000001490325489683642
$ws.Range("C6").Value = -0.000006155518350703382
$ws.Range("D6").Value = -0.001213179381647222
$ws.Range("F6").Value = -0.000006155518350703382
$ws.Range("G6").Value = -0.00003077759175351691
$ws.Range("C7").Value = -0.00000006381647378594923
$ws.Range("D7").Value = -0.0000002552658951437969
$ws.Range("E7").Value = -0.0002918794942274872
$ws.Range("F7").Value = -0.00004720066380059507
$ws.Range("G7").Value = -0.0000003190823689158684
$ws.Range("B8").Value = -0.001273030135507724
$ws.Range("C8").Value = -0.000001899458726306591
$ws.Range("D8").Value = -0.000176781043705887
$ws.Range("F8").Value = -0.000001899458726306591
$ws.Range("G8").Value = -0.000009497293632421133
$ws.Range("C9").Value = 0.0000008617006481692613
$ws.Range("D9").Value = 0.0007424967554925388
$ws.Range("F9").Value = 0.0000008617006481692613
$ws.Range("G9").Value = 0.000004308503241290396
$ws.Range("C10").Value = -0.0003111883291566642
$ws.Range("D10").Value = -0.001246275159303423
$ws.Range("E10").Value = -0.002983686692886067
$ws.Range("F10").Value = -0.00239223506477515
$ws.Range("G10").Value = -0.001555941645790426
$ws.Range("C11").Value = 0.000005793568329082177
$ws.Range("D11").Value = -0.0004955040387812915
$ws.Range("F11").Value = -0.0001359001812488714
$ws.Range("G11").Value = -0.01078138956654584
